$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.515.10"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.624.27"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'211.84"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'23.30"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "'0.0611"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "1.854.21"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.616.23"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'0.552"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "'65.43"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "27.490.20"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'229.88"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +8.51%  "
$ws.Range("D25").Value = "'149.16"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'15.52"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "'0.0484"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D33").Value = "1.468.30"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'3.06"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "'0.945"
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "'0.873"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.554"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'67.80"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'5.33"
$ws.Range("E46").Value = "  -4.89%  "
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("D48").Value = "1.764.02"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'87.27"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  +0.53%  "
Write-Host "Updated crypto prices and volumes"
